$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header / summary block text + numbers (Cant. Trabajadores box area)
# ---------------------------------------------------------------------------
$ws.Range("D2").Value  = "ESTADO DE CUENTA"
$ws.Range("B7").Value  = "RAZON SOCIAL:"
$ws.Range("B11").Value = "VALOR MORA"
$ws.Range("E11").Value = 1526322
$ws.Range("B13").Value = "Cant. Trabajadores"
$ws.Range("E13").Value = "Cant. Periodos"
$ws.Range("F13").Value = 42
$ws.Range("H15").Value = "Novedad de Ingreso"
$ws.Range("I15").Value = "Novedad de Retiro"
$ws.Range("J15").Value = "Observaciones"

# ---------------------------------------------------------------------------
# 2. Refresh the "Periodo Mora" column (rows 16-57) with the new, ascending
#    period list (oldest -> newest), replacing the previous descending list.
# ---------------------------------------------------------------------------
$periods = @(
  "2203","2204","2205","2206","2207","2208","2209","2210","2211","2212",
  "2301","2302","2303","2304","2305","2306","2307","2308","2309","2310","2311","2312",
  "2401","2402","2403","2404","2405","2406","2407","2408","2409","2410","2411","2412",
  "2501","2502","2503","2504","2505","2506","2507","2508"
)

$row = 16
foreach ($p in $periods) {
  $ws.Range("E$row").Value = $p
  $row = $row + 1
}

# ---------------------------------------------------------------------------
# 3. Row 57 becomes the new last data row of the table, so it must carry the
#    heavier "closing" bottom border that used to belong to row 63.
# ---------------------------------------------------------------------------
$lastRow = $ws.Range("B57:J57")
$lastRow.Borders.Item(9).LineStyle = 1
$lastRow.Borders.Item(9).Weight = 2
$lastRow.Borders.Item(9).ColorIndex = 1

# ---------------------------------------------------------------------------
# 4. Drop the now-obsolete rows 58-63 (old periods 2202 down to 2108); this
#    shifts the signature-block rows (old 68/69) up to become rows 62/63.
# ---------------------------------------------------------------------------
$ws.Range("A58:A63").EntireRow.Delete()
